$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsCCM   = $wb.Worksheets.Item("CSC-CSCCCMvSoECBtY")
$wsSoC   = $wb.Worksheets.Item("CSC-CSCSoCECBiaSY")

# --- CSC-CSCCCMvSoECBtY: rebuild the "cost multiplier" axis (row 2) and the
# derived "share of existing capacity" axis (row 1) so it better represents
# diminishing marginal returns. ---

# Row 2 (cost multiplier) new values for C2:N2
$wsCCM.Range("C2").Value = 1.25
$wsCCM.Range("D2").Value = 1.5
$wsCCM.Range("E2").Value = 1.75
$wsCCM.Range("F2").Value = 2
$wsCCM.Range("G2").Value = 2.25
$wsCCM.Range("H2").Value = 2.5
$wsCCM.Range("I2").Value = 2.75
$wsCCM.Range("J2").Value = 3
$wsCCM.Range("K2").Value = 3.25
$wsCCM.Range("L2").Value = 3.5
$wsCCM.Range("M2").Value = 3.75
$wsCCM.Range("N2").Value = 4

# Row 1 (share of existing capacity) now computed from row 2 via formula
$wsCCM.Range("C1").Formula = "=0.25^(1/B2)"
$wsCCM.Range("D1:O1").Formula = "=0.25^(1/C2)"

# --- CSC-CSCSoCECBiaSY: bump "share of cost effective capacity built in a
# single year" from 0.4 to 0.75 for all technology rows that used 0.4
# (rows 2:15 and 18:25), leaving rows 16:17 untouched. ---
$wsSoC.Range("B2:AE15").Value = 0.75
$wsSoC.Range("B18:AE25").Value = 0.75

# --- Restore/update the view selections recorded in each sheet. The
# "About" sheet must stay the active tab. ---
$wsCCM.Range("A4:P14").Select() | Out-Null
$wsAbout.Range("AU56").Select() | Out-Null
